# Updates on specific rates calculation with the values from metabolite analyzer
#
# For the two "Avg_*" sheets, insert 4 new metabolite-analyzer columns
# (Glucose (Met), Glutamate (Met), Glutamine (Met), Lactate (Met)) right
# before the existing "NH4" column, rename that column to "NH4 (Met)",
# and populate the new columns with the data supplied by the analyzer.

$wb = $excel.ActiveWorkbook

# New data (rows 2..7) for the 4 inserted columns C,D,E,F on each sheet.
$sheet2Data = @{
    2 = @(-1.076030475392915, -0.02814899053539199, -0.3721631139368082, 2.070590951829618)
    3 = @(-0.2691266020995701, 0.005055772999066666, 0, 0.1985270286203769)
    4 = @(-0.5568743780993389, 0.01607931349766054, -0.09455813630085891, 0.7089885652960982)
    5 = @(-0.3112096859122847, -0.002687723818444813, -0.01426653739615039, 0.1170111569730567)
    6 = @(-0.2552810321423932, 0.001958795284982216, 0, 0.112366412596784)
    7 = @(-0.6593168666857534, 0.006892731005801132, 0, 0.199177560244001)
}

$sheet3Data = @{
    2 = @(-0.5082603650054142, -0.04289242023381354, -0.3462310745509806, 0)
    3 = @(-0.1360752038108523, 0.00574595555923754, 0, 0)
    4 = @(-0.3532059231582924, 0.01290609316766033, -0.1017071763425346, 0)
    5 = @(-0.2182484002289638, 0.001109751869038702, -0.01557784900746719, 0)
    6 = @(-0.1739291432736104, 0.0019606893911018, -0.002034388085680979, 0)
    7 = @(-0.3490911028167712, 0.002422819492199735, 0, 0)
}

$sheetsToUpdate = @(
    @{ Name = "Avg_U1_U3"; Data = $sheet2Data },
    @{ Name = "Avg_U4_U8"; Data = $sheet3Data }
)

foreach ($entry in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # Insert 4 new blank columns before the current column C ("NH4").
    $ws.Range("C1:F1").EntireColumn.Insert()

    # New header labels for the inserted columns.
    $ws.Range("C1").Value = "Glucose (Met)"
    $ws.Range("D1").Value = "Glutamate (Met)"
    $ws.Range("E1").Value = "Glutamine (Met)"
    $ws.Range("F1").Value = "Lactate (Met)"

    # The old "NH4" column (now shifted to G) becomes "NH4 (Met)".
    $ws.Range("G1").Value = "NH4 (Met)"

    # Populate the new metabolite-analyzer data columns for rows 2..7.
    $data = $entry.Data
    foreach ($row in 2..7) {
        $vals = $data[$row]
        $ws.Cells.Item($row, 3).Value = $vals[0]
        $ws.Cells.Item($row, 4).Value = $vals[1]
        $ws.Cells.Item($row, 5).Value = $vals[2]
        $ws.Cells.Item($row, 6).Value = $vals[3]
    }
}
